$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Copy the formatting pattern of row 19 (B=money-style, C/D=plain, F=text)
# onto row 24 before writing values, so the new row inherits the same
# per-cell styling as the other "no-CV" rows (19-22) instead of Excel's
# default column style.
$ws.Range("A19:F19").Copy() | Out-Null
$ws.Range("A24").PasteSpecial(-4122) | Out-Null

$ws.Range("A24").Value = "emotionSignificant"
$ws.Range("B24").Value = 0.94631730000000003
$ws.Range("C24").Value = 0.93213999999999997
$ws.Range("D24").Value = 0.90183999999999997
$ws.Range("F24").Value = "增加了dayCount，并且用glm验证了significant"

# The copy/paste also stamped an (empty) E24 - remove it so the row matches
# the intended A:D + F layout with no cell in E.
$ws.Range("E24").Clear() | Out-Null

$ws.Range("C24").Select() | Out-Null
